# Fix sale quantity prediction
# Amend result into integral value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Truncate predicted sale quantities (column C) for rows 2241-2515
# (the forecast tail of the series) down to integer values.
$rng = $ws.Range("C2241:C2515")
foreach ($cell in $rng.Cells) {
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = [int]$v
    }
}

# Reposition the saved view/selection the way the author left it:
# scrolled so row 2219 is at the top, with the amended column selected.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 2219
$win.ScrollColumn = 1

$ws.Range("C2241:C2515").Select() | Out-Null
